$d = $word.ActiveDocument

$replacements = @(
    @("2024-03-29 Friday", "2024-03-30 Saturday"),
    @("200÷2=", "891÷6="),
    @("567÷8=", "131÷2="),
    @("496÷5=", "728÷3="),
    @("948÷3=", "125÷7="),
    @("723÷5=", "392÷8="),
    @("749÷3=", "659÷7="),
    @("560÷2=", "416÷8="),
    @("957÷3=", "270÷4="),
    @("982÷4=", "891÷4="),
    @("266÷6=", "443÷4="),
    @("662÷4=", "164÷3="),
    @("375÷5=", "106÷8="),
    @("545÷2=", "500÷9="),
    @("984÷8=", "797÷4="),
    @("350÷8=", "634÷6="),
    @("739÷4=", "185÷3="),
    @("374÷8=", "919÷5="),
    @("930÷2=", "985÷5="),
    @("687÷3=", "705÷5="),
    @("847÷2=", "907÷9="),
    @("651÷3=", "588÷8="),
    @("218÷9=", "636÷5="),
    @("949÷4=", "389÷9="),
    @("924÷3=", "274÷9="),
    @("985÷4=", "931÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
